$wb = $excel.ActiveWorkbook

# Sheets
$wsWeek2 = $wb.Worksheets.Item("Nädal 2")
$wsWeek1 = $wb.Worksheets.Item("Nädal 1")

# --- "Nädal 2" sheet updates ---
# Row 15: Stop time (D15, 14:00) and Delta Time in minutes (F15)
$wsWeek2.Range("D15").Value = 0.58333333333333337
$wsWeek2.Range("F15").Value = 160

# Row 16: Date (B16, 09.02.2020) and Start time (C16, 09:17), Activity (G16)
$wsWeek2.Range("B16").Value = 43870
$wsWeek2.Range("C16").Value = 0.38680555555555557
$wsWeek2.Range("G16").Value = "Kodutöö 2"

# Total row: add "minutes" label next to the sum
$wsWeek2.Range("G19").Value = "minutes"

# --- "Nädal 1" sheet updates ---
$wsWeek1.Range("G19").Value = "minutes"

# Update the active selection on this sheet as well
$wsWeek1.Range("G19").Select()

# Re-activate "Nädal 2" and restore its selection, keeping it the visible tab
$wsWeek2.Activate()
$wsWeek2.Range("G19").Select()
